# Auto-generated edit applying the "Updated cryptos list" diff.
# Strategy:
#  - Column D (Price) values are digit/dot strings that Excel's COM
#    Range.Value coercion would otherwise misparse as numbers (dropping
#    trailing zeros, switching to scientific notation, mangling the
#    subscript-digit glyphs used for very small prices). Prefixing with
#    a literal apostrophe forces literal-text entry (like typing '0.492
#    into Excel by hand); Style is then reset to "Normal" so the cell
#    does not pick up a stray quote-prefix / text-number-format style
#    that the apostrophe entry mode would otherwise leave behind.
#  - Columns B, C, E are plain text / whitespace-padded percentages that
#    Excel stores as text natively, so a direct .Value assignment is fine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.540.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.22%  '
$ws.Range("D3").Value = '''1.575.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  +0.89%  '
$ws.Range("D5").Value = '''212.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").Value = '''0.492'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("E7").Value = '  +0.83%  '
$ws.Range("D8").Value = '''46.17'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.61%  '
$ws.Range("D9").Value = '''24.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.21%  '
$ws.Range("E10").Value = '  -0.86%  '
$ws.Range("D11").Value = '''0.0592'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("D12").Value = '''0.0880'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").Value = '''1.800.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").Value = '''1.578.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.523'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '''3.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("D17").Value = '''28.516.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.17%  '
$ws.Range("D18").Value = '''62.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("D19").Value = '''229.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '''7.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").Value = '''0.0' + [char]0x2083 + '0693'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("E23").Value = '  -4.40%  '
$ws.Range("D24").Value = '''9.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").Value = '''2.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.87%  '
$ws.Range("D26").Value = '''151.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = '''15.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("D28").Value = '''6.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("E29").Value = '  -2.14%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("D31").Value = '''1.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.84%  '
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("D33").Value = '''3.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.36%  '
$ws.Range("D34").Value = '''3.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.48%  '
$ws.Range("D35").Value = '''1.392.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("D36").Value = '''1.55'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.60%  '
$ws.Range("E37").Value = '  -3.38%  '
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.74%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '''0.0166'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.83%  '
$ws.Range("D41").Value = '''0.532'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.73%  '
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("D43").Value = '''0.794'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.48%  '
$ws.Range("D44").Value = '''5.63'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.11%  '
$ws.Range("D45").Value = '''1.85'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.11%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").Value = '''62.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.96%  '
$ws.Range("D48").Value = '''1.712.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").Value = '''85.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("D50").Value = '''0.0' + [char]0x2086 + '0103'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.09%  '
$ws.Range("D51").Value = '''0.0518'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.41%  '
